$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New expense row (row 7): payment registered for Valdeci covering
# Rogerio's vacation, paid 2026-02-03.
$ws.Range("B7").Value = "síndica"
$ws.Cells.Item(7, 3).Value = "2026-02-03"
$ws.Range("D7").Value = "Valdeci"
$ws.Range("E7").Value = 350
$ws.Range("E7").NumberFormat = '"$"#,##0_);[Red]("$"#,##0)'
$ws.Range("F7").Value = "cobrir férias Rogerio"
$ws.Range("G7").Value = "260203_Comprovante_Pgto_Valdeci_Servicos_Jan_2026"
$ws.Range("I7").Value = "sim"

# Widen the "nf" column so the longer new entry is readable.
$ws.Range("G1:G600").ColumnWidth = 32.65

# Restore the cursor/selection to just below the new row, as left by the author.
$null = $ws.Range("C8").Select()
